# Add the new temperature sensor descriptions to column B of Sheet1
# (rows 26-29, next to pin numbers 24-27).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B26").Value = "HLT Temp Sensor"
$ws.Range("B27").Value = "MT Temp Sensor"
$ws.Range("B28").Value = "BP Temp Sensor"
$ws.Range("B29").Value = "Fermenter Temp Sensor"

# Match the author's scrolled/selected view state: scrolled down so row 2
# is at the top, with B30 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
